$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 892, shifting the rest of the
# table (old rows 892-909) down to 895-912.
$ws.Rows("892:894").Insert()

# Row 892
$ws.Range("A892").Value = 3
$ws.Range("B892").Value = "Femacal de La Calera"
$ws.Range("C892").Value = "Coquimbo"
$ws.Range("D892").Value = 44448
$ws.Range("E892").Value = 5
$ws.Range("F892").Value = 100112002
$ws.Range("G892").Value = "Pimiento"
$ws.Range("H892").Value = "Zafiro rojo"
$ws.Range("I892").Value = "Segunda"
$ws.Range("J892").Value = 35
$ws.Range("K892").Value = 36000
$ws.Range("L892").Value = 36000
$ws.Range("M892").Value = 36000
$ws.Range("N892").Value = "`$/caja 15 kilos"
$ws.Range("O892").Value = "Región de Arica y Parinacota"
$ws.Range("P892").Value = 2400
$ws.Range("Q892").Value = 15
$ws.Range("R892").Value = "Hortaliza"

# Row 893
$ws.Range("A893").Value = 3
$ws.Range("B893").Value = "Femacal de La Calera"
$ws.Range("C893").Value = "Coquimbo"
$ws.Range("D893").Value = 44448
$ws.Range("E893").Value = 5
$ws.Range("F893").Value = 100112002
$ws.Range("G893").Value = "Pimiento"
$ws.Range("H893").Value = "Zafiro verde"
$ws.Range("I893").Value = "Primera"
$ws.Range("J893").Value = 73
$ws.Range("K893").Value = 36000
$ws.Range("L893").Value = 37000
$ws.Range("M893").Value = 36479
$ws.Range("N893").Value = "`$/caja 15 kilos"
$ws.Range("O893").Value = "Región de Arica y Parinacota"
$ws.Range("P893").Value = 2432
$ws.Range("Q893").Value = 15
$ws.Range("R893").Value = "Hortaliza"

# Row 894
$ws.Range("A894").Value = 3
$ws.Range("B894").Value = "Femacal de La Calera"
$ws.Range("C894").Value = "Coquimbo"
$ws.Range("D894").Value = 44448
$ws.Range("E894").Value = 5
$ws.Range("F894").Value = 100112002
$ws.Range("G894").Value = "Pimiento"
$ws.Range("H894").Value = "Zafiro verde"
$ws.Range("I894").Value = "Segunda"
$ws.Range("J894").Value = 38
$ws.Range("K894").Value = 32000
$ws.Range("L894").Value = 32000
$ws.Range("M894").Value = 32000
$ws.Range("N894").Value = "`$/caja 15 kilos"
$ws.Range("O894").Value = "Región de Arica y Parinacota"
$ws.Range("P894").Value = 2133
$ws.Range("Q894").Value = 15
$ws.Range("R894").Value = "Hortaliza"
